$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Global rename of the "short-url" shared string used by every data row
# (column B, B2:B580). "PR4eI5" is not numeric-looking, so Excel will
# naturally store it as text using the existing left-aligned style.
$ws.Range("B2:B580").Value = "PR4eI5"

# Row 562 - Afghanistan -> Pakistan, year 2024: revised refugee figures.
# These look like numbers, so a leading apostrophe keeps them as text
# (matching the source data's convention of storing numbers as strings).
$ws.Range("N562").Value = "'1559964"
$ws.Range("O562").Value = "'198610"
$ws.Range("P562").Value = "'126812"
$ws.Range("T562").Value = "'798348"

# Row 565 - China: asylum_seekers
$ws.Range("O565").Value = "'5"

# Row 567 - Iran (Islamic Rep. of): asylum_seekers
$ws.Range("O567").Value = "'32"

# Row 570 - Myanmar: refugees
$ws.Range("N570").Value = "'28"

# Row 571 - Pakistan -> Pakistan, year 2024: idps / ooc / hst
$ws.Range("Q571").Value = "'3439"
$ws.Range("T571").Value = "'283"
$ws.Range("V571").Value = "'1456000"

# Row 573 - Somalia: asylum_seekers
$ws.Range("O573").Value = "'34"

# Row 574 - Stateless: stateless
$ws.Range("S574").Value = "'60"

# Row 575 - Sudan: asylum_seekers
$ws.Range("O575").Value = "'11"

# Row 578 - Tajikistan: asylum_seekers
$ws.Range("O578").Value = "'8"

# Row 580 - Yemen: asylum_seekers
$ws.Range("O580").Value = "'44"
